# Added date conversion method
# Copy the "Date of Birth" / "Date of Vaccination" columns from the "aus"
# sheet into the "ind" sheet as new columns C/D, and update the selections
# on both sheets.

$wb = $excel.ActiveWorkbook

# --- "ind" sheet: add Date of Birth / Date of Vaccination columns -------
$ws = $wb.Worksheets.Item("ind")

# Header row
$ws.Range("C1").Value = "Date of Birth"
$ws.Range("D1").Value = "Date of Vaccination"

# Row 2 (Mike) - birth date unknown / vaccination date known
$ws.Range("C2").Value = "NULL"
$ws.Range("D2").Value = 44692
$ws.Range("D2").NumberFormat = "dd/mm/yyyy"

# Row 3 (Jonnathan) - birth date known / vaccination date invalid text
$ws.Range("C3").Value = 35777
$ws.Range("C3").NumberFormat = "dd/mm/yyyy"
$ws.Range("D3").Value = "2021-13-13"

# Row 4 (Cristina) - both dates known
$ws.Range("C4").Value = 35866
$ws.Range("C4").NumberFormat = "dd/mm/yyyy"
$ws.Range("D4").Value = 44632
$ws.Range("D4").NumberFormat = "dd/mm/yyyy"

# New column widths for the two inserted columns
$ws.Columns.Item(3).ColumnWidth = 11
$ws.Columns.Item(4).ColumnWidth = 16.5

# --- Update view selections to match the edited state -------------------
$aus = $wb.Worksheets.Item("aus")
$null = $aus.Range("B1").Select()

$null = $ws.Range("C7").Select()
